{"js": "async (context) => {\n  // Replacement text for the document's first 12 paragraphs (0-based indices 0-11),\n  // in document order. This mirrors the target OOXML from the diff, including the\n  // single entry (index 4) that keeps a trailing space (Word will automatically add\n  // xml:space=\"preserve\" when it saves a run whose text has leading/trailing spaces).\n  const newParagraphTexts = [\n    \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 05.09.25\",\n    \"Group Sequence Policy Optimization\",\n    \"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05de\u05e1' 500:\",\n    \"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05de\u05e1' 500 \u05d5\u05d4\u05d3\u05d9 \u05d7\u05d2\u05d9\u05d2\u05d9\u05ea \u05dc\u05db\u05d0\u05d5\u05e8\u05d4, \u05d1\u05d4\u05ea\u05d7\u05dc\u05d4 \u05d7\u05e9\u05d1\u05ea\u05d9 \u05dc\u05d1\u05d7\u05d5\u05e8 \u05d0\u05d9\u05d6\u05d4 \u05de\u05d0\u05de\u05e8 \u05de\u05d9\u05d5\u05d7\u05d3 \u05d0\u05d1\u05dc \u05dc\u05d0\u05d7\u05e8 \u05d4\u05e8\u05d4\u05d5\u05e8\u05d9\u05dd \u05e2\u05de\u05d5\u05e7\u05d9\u05dd (\u05d0\u05da \u05dc\u05d0 \u05d0\u05e8\u05d5\u05db\u05d9\u05dd) \u05d4\u05d7\u05dc\u05d8\u05ea\u05d9 \u05dc\u05d3\u05d7\u05d5\u05ea \u05d0\u05ea \u05d4\u05d7\u05d2\u05d9\u05d2\u05d4 \u05dc\u05de\u05d0\u05de\u05e8 \u05de\u05e1' 512. \u05d5\u05e9\u05dd \u05db\u05d1\u05e8 \u05e0\u05d7\u05dc\u05d9\u05d8, \u05d0\u05d5\u05dc\u05d9 \u05e0\u05d3\u05d7\u05d4 \u05dc-555 \u05d0\u05d5 \u05de\u05e9\u05d4\u05d5 \u05db\u05d6\u05d4 - \u05e0\u05e8\u05d0\u05d4 \u05d0\u05d9\u05da \u05d4\u05ea\u05e7\u05d3\u05de\u05d5 \u05d4\u05d4\u05e4\u05ea\u05e2\u05d5\u05ea \u05e9\u05e9\u05d5\u05ea\u05e4\u05d9\u05d9 \u05d5\u05d0\u05e0\u05d9 \u05de\u05db\u05d9\u05e0\u05d9\u05dd \u05dc\u05db\u05dd \ud83d\ude42.\",\n    \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05e9\u05db\u05dc\u05d5\u05dc \u05dc\u05e9\u05d9\u05d8\u05ea GRPO \u05d0\u05d5 Group Relative Policy Optimization \u05d4\u05e9\u05d9\u05d9\u05db\u05ea \u05dc\u05de\u05e9\u05e4\u05d7\u05ea \u05e9\u05d9\u05d8\u05d5\u05ea RHLF \u05d4\u05de\u05e9\u05de\u05e9\u05d5\u05ea \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05d5\u05dc\u05e4\u05d9\u05d9\u05e0\u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e9\u05e7\u05d9\u05d1\u05dc\u05d4 \u05e9\u05dd GSPO (\u05d4\u05d7\u05dc\u05d9\u05e4\u05d5 Relative \u05d1- Sequence) \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05e9\u05e0\u05d4 \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc GRPO. \",\n    \"\u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 GRPO \u05de\u05de\u05e7\u05e1\u05dd \u05d0\u05ea \u05d4\u05de\u05db\u05e4\u05dc\u05d4 \u05e9\u05dc \u05e9\u05e0\u05d9 \u05d4\u05d2\u05d5\u05e8\u05de\u05d9\u05dd(\u05d9\u05e9 \u05d2\u05dd \u05db\u05de\u05d4 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e7\u05dc\u05d9\u05e4 \u05e9\u05dd). \u05d4\u05d2\u05d5\u05e8\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 \u05de\u05d4 \u05d4\u05d5\u05d0 \u05d4\u05d9\u05ea\u05e8\u05d5\u05df \u05e9\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9 (\u05e9\u05d6\u05d4 \u05d1\u05e2\u05e6\u05dd \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05ea\u05e0\u05d9\u05ea \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05d4\u05e7\u05e9\u05e8 \u05d4\u05e7\u05d5\u05d3\u05dd \u05dc\u05d5) \u05e2\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d9\u05e9\u05df (\u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d2\u05d3\u05de\u05d9\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df). GRPO \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de-PPO \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9 \u05dc\u05d0 \u05de\u05d7\u05e9\u05d1 \u05d0\u05d5\u05ea\u05d4 \u05d3\u05e8\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea value \u05d0\u05dc\u05d0 \u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd (rewards) \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d4\u05de\u05ea\u05e7\u05d1\u05dc\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05e0\u05d3\u05d2\u05de\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 (\u05d1\u05d2\u05dc\u05dc \u05d6\u05d4 \u05de\u05d9\u05dc\u05d4 group \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05d1\u05e9\u05dd \u05e9\u05dc \u05d4\u05e9\u05d9\u05d8\u05d4).\",\n    \"\u05d4\u05d2\u05d5\u05e8\u05dd \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05d9\u05d7\u05e1 \u05e9\u05dc \u05e4\u05d5\u05dc\u05d5\u05e1\u05d9 \u05d4\u05d7\u05d3\u05e9 \u05e9\u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d0\u05e4\u05d8\u05de\u05d9\u05dd (\u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05ea\u05e0\u05d9\u05ea \u05e9\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc) \u05dc\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d9\u05e9\u05df \u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d2\u05d3\u05de\u05d9\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd. \u05db\u05d0\u05df \u05d1\u05d0 \u05d4\u05d4\u05d1\u05d3\u05dc \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9 \u05d1\u05d9\u05df GRPO \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e0\u05de\u05e6\u05d0 \u05d1\u05d0\u05d9\u05da \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05d9\u05d7\u05e1 \u05d4\u05d6\u05d4. \u05d1-GRPO \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d7\u05d3\u05e9 \u05d5\u05d4\u05d9\u05e9\u05df \u05d1\u05e8\u05de\u05ea \u05d4\u05d8\u05d5\u05e7\u05df \u05de\u05e0\u05d5\u05e8\u05de\u05dc\u05d9 \u05d1\u05d0\u05d5\u05e8\u05da \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05d3 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d6\u05d4. \u05d7\u05d9\u05e9\u05d5\u05d1 \u05d6\u05d4 \u05db\u05de\u05d5\u05d1\u05df \u05d1\u05e2\u05dc \u05e9\u05d5\u05e0\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4 \u05d5\u05d6\u05d4 \u05d4\u05e1\u05d9\u05d1\u05d4 \u05dc\u05d4\u05d9\u05de\u05e6\u05d0\u05d5\u05ea \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dd \u05db\u05de\u05d4 \u05e7\u05dc\u05d9\u05e4\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05de\u05e0\u05d5\u05e2 \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8. \u05d3\u05f4\u05d0 \u05d1-PPO \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05e8\u05de\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05d5\u05dc\u05d4 \u05d0\u05d1\u05dc \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05dc\u05d3\u05dc\u05d9\u05dc\u05d9\u05dd (sparse) \u05e9\u05d6\u05d4 \u05db\u05de\u05d5\u05d1\u05df \u05ea\u05e8\u05d7\u05d9\u05e9 \u05dc\u05d0 \u05e4\u05e9\u05d5\u05d8 \u05d1\u05d1\u05e2\u05d9\u05d5\u05ea RL.\",\n    \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 2 \u05e9\u05d9\u05d8\u05d5\u05ea. \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4, \u05e9\u05de\u05d7\u05d6\u05d9\u05e8\u05d4 \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d1\u05e8\u05de\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05d5\u05dc\u05d4, \u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05ea \u05d4\u05d9\u05d7\u05e1 \u05d1\u05ea\u05d5\u05e8 \u05de\u05de\u05d5\u05e6\u05e2 \u05e2\u05dc \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd (\u05d1\u05dc\u05d5\u05d2 \u05e1\u05e7\u05d9\u05d9\u05dc) \u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05d0\u05d7\u05ea \u05de\u05d4\u05df \u05de\u05e0\u05d5\u05e8\u05de\u05dc\u05ea \u05d1\u05d0\u05d5\u05e8\u05da \u05e9\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05d3 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d6\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05d4\u05e9\u05e0\u05d9\u05d4 \u05de\u05e9\u05d0\u05d9\u05e8\u05d4, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-GRPO, \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d1\u05e8\u05de\u05ea \u05d4\u05d8\u05d5\u05e7\u05df \u05d0\u05d1\u05dc \u05d9\u05d7\u05e1 \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05de\u05d7\u05d5\u05e9\u05d1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 - \u05e8\u05e7 \u05e9\u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05de\u05d7\u05d5\u05e9\u05d1 \u05e2\u05dc \u05d4\u05d8\u05d5\u05e7\u05df. \u05e9\u05ea\u05d9 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea  \u05e0\u05e8\u05d0\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de- GRPO \u05d0\u05d1\u05dc \u05d4\u05e7\u05dc\u05d9\u05e4\u05d9\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e0\u05de\u05e6\u05d0\u05d9\u05dd \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05d4 \u05d4\u05de\u05d8\u05e8\u05d4.\",\n    \"\u05d9\u05e9 \u05d1\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d8\u05e2\u05e0\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05d1\u05d9\u05df \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05d5\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc GRPO \u05d5-PPO \u05dc-importance sampling \u05d0\u05d5 IM. \u05d0\u05d6\u05db\u05d9\u05e8 \u05db\u05d9 IM \u05d4\u05d9\u05d0 \u05e9\u05d9\u05d8\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea P \u05e9\u05e7\u05e9\u05d4 \u05dc\u05d3\u05d2\u05d5\u05dd \u05de\u05de\u05e0\u05d5 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea Q \u05e9\u05d9\u05d5\u05ea\u05e8 \u05e7\u05dc \u05dc\u05d3\u05d2\u05d5\u05dd \u05de\u05de\u05e0\u05d4. \u05de\u05e9\u05e7\u05dc \u05d4-importance \u05e2\u05d1\u05d5\u05e8 \u05d3\u05d2\u05d9\u05de\u05d4 x \u05d4\u05d5\u05d0 \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05ea \u05e9\u05dc x \u05e2\u05dd P \u05d5\u05e2\u05dd Q. \u05d0\u05de\u05e0\u05dd \u05d9\u05e9 \u05e7\u05e9\u05e8 \u05d0\u05de\u05d9\u05ea\u05d9 \u05d1\u05d9\u05df IM \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05de\u05d5\u05d6\u05db\u05e8\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05e9\u05ea\u05db\u05e0\u05e2\u05ea\u05d9 \u05e9\u05db\u05dc \u05d4\u05e0\u05d9\u05de\u05d5\u05e7\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05dd \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05de\u05ea\u05de\u05d8\u05d9\u05ea - \u05d9\u05e9 \u05de\u05e6\u05d1 \u05e9\u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05d1\u05e0\u05ea\u05d9 \u05d0\u05d5\u05ea\u05dd \u05de\u05e1\u05e4\u05d9\u05e7 \u05e2\u05de\u05d5\u05e7.\",\n    \"\u05db\u05da \u05d0\u05d5 \u05db\u05da \u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05e8\u05d0\u05d5\u05d9 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05e4\u05e8 500!\",\n    \"https://www.arxiv.org/abs/2507.18071\",\n    \"\u05d4\u05d0\u05d5\u05ea \u05d5\u05d4\u05e8\u05e2\u05e9: \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc\u05d4\u05e2\u05e8\u05db\u05ea \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4\"\n  ];\n\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  // Overwrite the text of the first 12 paragraphs in place (this keeps each\n  // paragraph's own <w:p>/<w:pPr> and produces a single clean <w:r><w:t>... run,\n  // matching the diff).\n  for (let i = 0; i < newParagraphTexts.length; i++) {\n    paragraphs.items[i].insertText(newParagraphTexts[i], Word.InsertLocation.replace);\n  }\n\n  // The remaining trailing paragraphs (originally indices 12-16, i.e. the\n  // \"\u05d1\u05d7\u05d9\u05e8\u05ea \u05d7\u05dc\u05e7\u05d9\u05dd...\" / \"\u05d4\u05e4\u05d7\u05ea\u05ea \u05e8\u05e2\u05e9...\" / \"\u05e9\u05d9\u05e0\u05d5\u05d9 \u05e1\u05d5\u05d2 \u05d4\u05de\u05d3\u05d9\u05d3\u05d4...\" / \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e1\u05d9\u05e4\u05e7\u05d5...\" /\n  // \"https://arxiv.org/abs/2508.13144\" paragraphs) are removed entirely, from the\n  // end backwards so the indices of paragraphs still to be deleted stay valid.\n  for (let i = paragraphs.items.length - 1; i >= newParagraphTexts.length; i--) {\n    paragraphs.items[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Grab the active document (already open in the headless Word COM host).\n$d = $word.ActiveDocument\n\n# Replacement text for the document's first 12 paragraphs (1-based Word COM\n# indices 1-12), in document order. This mirrors the target OOXML from the diff.\n$newParagraphTexts = @(\n  '\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 05.09.25',\n  'Group Sequence Policy Optimization',\n  '\u05e1\u05e7\u05d9\u05e8\u05d4 \u05de\u05e1'' 500:',\n  '\u05e1\u05e7\u05d9\u05e8\u05d4 \u05de\u05e1'' 500 \u05d5\u05d4\u05d3\u05d9 \u05d7\u05d2\u05d9\u05d2\u05d9\u05ea \u05dc\u05db\u05d0\u05d5\u05e8\u05d4, \u05d1\u05d4\u05ea\u05d7\u05dc\u05d4 \u05d7\u05e9\u05d1\u05ea\u05d9 \u05dc\u05d1\u05d7\u05d5\u05e8 \u05d0\u05d9\u05d6\u05d4 \u05de\u05d0\u05de\u05e8 \u05de\u05d9\u05d5\u05d7\u05d3 \u05d0\u05d1\u05dc \u05dc\u05d0\u05d7\u05e8 \u05d4\u05e8\u05d4\u05d5\u05e8\u05d9\u05dd \u05e2\u05de\u05d5\u05e7\u05d9\u05dd (\u05d0\u05da \u05dc\u05d0 \u05d0\u05e8\u05d5\u05db\u05d9\u05dd) \u05d4\u05d7\u05dc\u05d8\u05ea\u05d9 \u05dc\u05d3\u05d7\u05d5\u05ea \u05d0\u05ea \u05d4\u05d7\u05d2\u05d9\u05d2\u05d4 \u05dc\u05de\u05d0\u05de\u05e8 \u05de\u05e1'' 512. \u05d5\u05e9\u05dd \u05db\u05d1\u05e8 \u05e0\u05d7\u05dc\u05d9\u05d8, \u05d0\u05d5\u05dc\u05d9 \u05e0\u05d3\u05d7\u05d4 \u05dc-555 \u05d0\u05d5 \u05de\u05e9\u05d4\u05d5 \u05db\u05d6\u05d4 - \u05e0\u05e8\u05d0\u05d4 \u05d0\u05d9\u05da \u05d4\u05ea\u05e7\u05d3\u05de\u05d5 \u05d4\u05d4\u05e4\u05ea\u05e2\u05d5\u05ea \u05e9\u05e9\u05d5\u05ea\u05e4\u05d9\u05d9 \u05d5\u05d0\u05e0\u05d9 \u05de\u05db\u05d9\u05e0\u05d9\u05dd \u05dc\u05db\u05dd \ud83d\ude42.',\n  '\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05e9\u05db\u05dc\u05d5\u05dc \u05dc\u05e9\u05d9\u05d8\u05ea GRPO \u05d0\u05d5 Group Relative Policy Optimization \u05d4\u05e9\u05d9\u05d9\u05db\u05ea \u05dc\u05de\u05e9\u05e4\u05d7\u05ea \u05e9\u05d9\u05d8\u05d5\u05ea RHLF \u05d4\u05de\u05e9\u05de\u05e9\u05d5\u05ea \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05d5\u05dc\u05e4\u05d9\u05d9\u05e0\u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e9\u05e7\u05d9\u05d1\u05dc\u05d4 \u05e9\u05dd GSPO (\u05d4\u05d7\u05dc\u05d9\u05e4\u05d5 Relative \u05d1- Sequence) \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05e9\u05e0\u05d4 \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc GRPO. ',\n  '\u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 GRPO \u05de\u05de\u05e7\u05e1\u05dd \u05d0\u05ea \u05d4\u05de\u05db\u05e4\u05dc\u05d4 \u05e9\u05dc \u05e9\u05e0\u05d9 \u05d4\u05d2\u05d5\u05e8\u05de\u05d9\u05dd(\u05d9\u05e9 \u05d2\u05dd \u05db\u05de\u05d4 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e7\u05dc\u05d9\u05e4 \u05e9\u05dd). \u05d4\u05d2\u05d5\u05e8\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 \u05de\u05d4 \u05d4\u05d5\u05d0 \u05d4\u05d9\u05ea\u05e8\u05d5\u05df \u05e9\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9 (\u05e9\u05d6\u05d4 \u05d1\u05e2\u05e6\u05dd \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05ea\u05e0\u05d9\u05ea \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05d4\u05e7\u05e9\u05e8 \u05d4\u05e7\u05d5\u05d3\u05dd \u05dc\u05d5) \u05e2\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d9\u05e9\u05df (\u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d2\u05d3\u05de\u05d9\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df). GRPO \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de-PPO \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9 \u05dc\u05d0 \u05de\u05d7\u05e9\u05d1 \u05d0\u05d5\u05ea\u05d4 \u05d3\u05e8\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea value \u05d0\u05dc\u05d0 \u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd (rewards) \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d4\u05de\u05ea\u05e7\u05d1\u05dc\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05e0\u05d3\u05d2\u05de\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 (\u05d1\u05d2\u05dc\u05dc \u05d6\u05d4 \u05de\u05d9\u05dc\u05d4 group \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05d1\u05e9\u05dd \u05e9\u05dc \u05d4\u05e9\u05d9\u05d8\u05d4).',\n  '\u05d4\u05d2\u05d5\u05e8\u05dd \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05d9\u05d7\u05e1 \u05e9\u05dc \u05e4\u05d5\u05dc\u05d5\u05e1\u05d9 \u05d4\u05d7\u05d3\u05e9 \u05e9\u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d0\u05e4\u05d8\u05de\u05d9\u05dd (\u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05de\u05d5\u05ea\u05e0\u05d9\u05ea \u05e9\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc) \u05dc\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d9\u05e9\u05df \u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d2\u05d3\u05de\u05d9\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd. \u05db\u05d0\u05df \u05d1\u05d0 \u05d4\u05d4\u05d1\u05d3\u05dc \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9 \u05d1\u05d9\u05df GRPO \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e0\u05de\u05e6\u05d0 \u05d1\u05d0\u05d9\u05da \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05d9\u05d7\u05e1 \u05d4\u05d6\u05d4. \u05d1-GRPO \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d4\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05d7\u05d3\u05e9 \u05d5\u05d4\u05d9\u05e9\u05df \u05d1\u05e8\u05de\u05ea \u05d4\u05d8\u05d5\u05e7\u05df \u05de\u05e0\u05d5\u05e8\u05de\u05dc\u05d9 \u05d1\u05d0\u05d5\u05e8\u05da \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05d3 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d6\u05d4. \u05d7\u05d9\u05e9\u05d5\u05d1 \u05d6\u05d4 \u05db\u05de\u05d5\u05d1\u05df \u05d1\u05e2\u05dc \u05e9\u05d5\u05e0\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4 \u05d5\u05d6\u05d4 \u05d4\u05e1\u05d9\u05d1\u05d4 \u05dc\u05d4\u05d9\u05de\u05e6\u05d0\u05d5\u05ea \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dd \u05db\u05de\u05d4 \u05e7\u05dc\u05d9\u05e4\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05de\u05e0\u05d5\u05e2 \u05e9\u05d9\u05e0\u05d5\u05d9\u05d9\u05dd \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8. \u05d3\u05f4\u05d0 \u05d1-PPO \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05e8\u05de\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05d5\u05dc\u05d4 \u05d0\u05d1\u05dc \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05dc\u05d3\u05dc\u05d9\u05dc\u05d9\u05dd (sparse) \u05e9\u05d6\u05d4 \u05db\u05de\u05d5\u05d1\u05df \u05ea\u05e8\u05d7\u05d9\u05e9 \u05dc\u05d0 \u05e4\u05e9\u05d5\u05d8 \u05d1\u05d1\u05e2\u05d9\u05d5\u05ea RL.',\n  '\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 2 \u05e9\u05d9\u05d8\u05d5\u05ea. \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4, \u05e9\u05de\u05d7\u05d6\u05d9\u05e8\u05d4 \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d1\u05e8\u05de\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05d5\u05dc\u05d4, \u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05ea \u05d4\u05d9\u05d7\u05e1 \u05d1\u05ea\u05d5\u05e8 \u05de\u05de\u05d5\u05e6\u05e2 \u05e2\u05dc \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd (\u05d1\u05dc\u05d5\u05d2 \u05e1\u05e7\u05d9\u05d9\u05dc) \u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05d0\u05d7\u05ea \u05de\u05d4\u05df \u05de\u05e0\u05d5\u05e8\u05de\u05dc\u05ea \u05d1\u05d0\u05d5\u05e8\u05da \u05e9\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05d3 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d6\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05d4\u05e9\u05e0\u05d9\u05d4 \u05de\u05e9\u05d0\u05d9\u05e8\u05d4, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-GRPO, \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d1\u05e8\u05de\u05ea \u05d4\u05d8\u05d5\u05e7\u05df \u05d0\u05d1\u05dc \u05d9\u05d7\u05e1 \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05ea \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05de\u05d7\u05d5\u05e9\u05d1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 - \u05e8\u05e7 \u05e9\u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05de\u05d7\u05d5\u05e9\u05d1 \u05e2\u05dc \u05d4\u05d8\u05d5\u05e7\u05df. \u05e9\u05ea\u05d9 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea  \u05e0\u05e8\u05d0\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de- GRPO \u05d0\u05d1\u05dc \u05d4\u05e7\u05dc\u05d9\u05e4\u05d9\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e0\u05de\u05e6\u05d0\u05d9\u05dd \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05d4 \u05d4\u05de\u05d8\u05e8\u05d4.',\n  '\u05d9\u05e9 \u05d1\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d8\u05e2\u05e0\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05d1\u05d9\u05df \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05d5\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc GRPO \u05d5-PPO \u05dc-importance sampling \u05d0\u05d5 IM. \u05d0\u05d6\u05db\u05d9\u05e8 \u05db\u05d9 IM \u05d4\u05d9\u05d0 \u05e9\u05d9\u05d8\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea P \u05e9\u05e7\u05e9\u05d4 \u05dc\u05d3\u05d2\u05d5\u05dd \u05de\u05de\u05e0\u05d5 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea Q \u05e9\u05d9\u05d5\u05ea\u05e8 \u05e7\u05dc \u05dc\u05d3\u05d2\u05d5\u05dd \u05de\u05de\u05e0\u05d4. \u05de\u05e9\u05e7\u05dc \u05d4-importance \u05e2\u05d1\u05d5\u05e8 \u05d3\u05d2\u05d9\u05de\u05d4 x \u05d4\u05d5\u05d0 \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d4\u05d4\u05e1\u05ea\u05d1\u05e8\u05d5\u05ea \u05e9\u05dc x \u05e2\u05dd P \u05d5\u05e2\u05dd Q. \u05d0\u05de\u05e0\u05dd \u05d9\u05e9 \u05e7\u05e9\u05e8 \u05d0\u05de\u05d9\u05ea\u05d9 \u05d1\u05d9\u05df IM \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05d4\u05de\u05d5\u05d6\u05db\u05e8\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05e9\u05ea\u05db\u05e0\u05e2\u05ea\u05d9 \u05e9\u05db\u05dc \u05d4\u05e0\u05d9\u05de\u05d5\u05e7\u05d9\u05dd \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05dd \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05de\u05ea\u05de\u05d8\u05d9\u05ea - \u05d9\u05e9 \u05de\u05e6\u05d1 \u05e9\u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05d1\u05e0\u05ea\u05d9 \u05d0\u05d5\u05ea\u05dd \u05de\u05e1\u05e4\u05d9\u05e7 \u05e2\u05de\u05d5\u05e7.',\n  '\u05db\u05da \u05d0\u05d5 \u05db\u05da \u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05e8\u05d0\u05d5\u05d9 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05e4\u05e8 500!',\n  'https://www.arxiv.org/abs/2507.18071',\n  '\u05d4\u05d0\u05d5\u05ea \u05d5\u05d4\u05e8\u05e2\u05e9: \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc\u05d4\u05e2\u05e8\u05db\u05ea \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4'\n)\n\n# Overwrite the text of each of those paragraphs in place. Assigning to\n# Range.Text replaces only the paragraph's content (not its paragraph mark),\n# so formatting/paragraph style and the <w:p>/<w:pPr> wrapper are preserved.\nfor ($i = 0; $i -lt $newParagraphTexts.Count; $i++) {\n  $d.Paragraphs.Item($i + 1).Range.Text = $newParagraphTexts[$i]\n}\n\n# The remaining trailing paragraphs (originally 13-17, i.e. the\n# \"\u05d1\u05d7\u05d9\u05e8\u05ea \u05d7\u05dc\u05e7\u05d9\u05dd...\" / \"\u05d4\u05e4\u05d7\u05ea\u05ea \u05e8\u05e2\u05e9...\" / \"\u05e9\u05d9\u05e0\u05d5\u05d9 \u05e1\u05d5\u05d2 \u05d4\u05de\u05d3\u05d9\u05d3\u05d4...\" / \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e1\u05d9\u05e4\u05e7\u05d5...\" /\n# \"https://arxiv.org/abs/2508.13144\" paragraphs) are deleted entirely in one go.\n$totalParagraphs = $d.Paragraphs.Count\nif ($totalParagraphs -gt $newParagraphTexts.Count) {\n  $firstParaToRemove = $newParagraphTexts.Count + 1\n  $deleteRange = $d.Range(\n    $d.Paragraphs.Item($firstParaToRemove).Range.Start,\n    $d.Paragraphs.Item($totalParagraphs).Range.End\n  )\n  $deleteRange.Delete()\n}\n"}
